$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before D. This shifts:
#    D(qty) -> E, G(formula) -> H, L,M,N (insert-string parts) -> M,N,O
$ws.Columns("D").Insert()

# 2. Insert the extra rows needed for the two split tickets (3283295 and 2374827)
#    Row numbers below are "as they exist at the moment of each call".
$ws.Rows(16).Insert()
$ws.Rows(18).Insert()
$ws.Rows(21).Insert()

# 3. Update the SQL literal text (now in M1) to the new "insert ... values(" text
$ws.Range("M1").Value = "insert into buysTicket values("

# 4. Fill in the new "date" column D with the purchase date for each ticket group,
#    and update the H-column formula to splice the date in before the closing ");"

# --- Group 1: rows 1-4, ticket 2321423, date 2016-03-20 (quote-prefixed, like the original file) ---
$ws.Range("D1:D4").Value = "' '2016-03-20'"
$ws.Range("H1").Formula = "=CONCATENATE(`$M`$1,A1,`$O`$1,B1,`$O`$1,C1,`$O`$1,D1,`$N`$1)"
$ws.Range("H2").Formula = "=CONCATENATE(`$M`$1,A2,`$O`$1,B2,`$O`$1,C2,`$O`$1,D2,`$N`$1)"
$ws.Range("H3").Formula = "=CONCATENATE(`$M`$1,A3,`$O`$1,B3,`$O`$1,C3,`$O`$1,D3,`$N`$1)"
$ws.Range("H4").Formula = "=CONCATENATE(`$M`$1,A4,`$O`$1,B4,`$O`$1,C4,`$O`$1,D4,`$N`$1)"

# --- Group 2: row 8, ticket 4737492, date 2015-02-04 ---
$ws.Range("D8").Value = " '2015-02-04'"
$ws.Range("H8").Formula = "=CONCATENATE(`$M`$1,A8,`$O`$1,B8,`$O`$1,C8,`$O`$1,D8,`$N`$1)"

# --- Group 3: row 10, ticket 1238471, date 2018-03-01 ---
$ws.Range("D10").Value = " '2018-03-01'"
$ws.Range("H10").Formula = "=CONCATENATE(`$M`$1,A10,`$O`$1,B10,`$O`$1,C10,`$O`$1,D10,`$N`$1)"

# --- Group 4: row 14, ticket 4917493, date 2013-12-12 ---
$ws.Range("D14").Value = " '2013-12-12'"
$ws.Range("H14").Formula = "=CONCATENATE(`$M`$1,A14,`$O`$1,B14,`$O`$1,C14,`$O`$1,D14,`$N`$1)"

# --- Group 5: rows 17 & 18, ticket 3283295, date 2009-07-29 (split across two seats) ---
$ws.Range("D17").Value = " '2009-07-29'"
$ws.Range("H17").Formula = "=CONCATENATE(`$M`$1,A17,`$O`$1,B17,`$O`$1,C17,`$O`$1,D17,`$N`$1)"

$ws.Range("A18").Value = 3283295
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 470
$ws.Range("D18").Value = " '2009-07-29'"
$ws.Range("H18").Formula = "=CONCATENATE(`$M`$1,A18,`$O`$1,B18,`$O`$1,C18,`$O`$1,D18,`$N`$1)"

# --- Group 6: rows 20 & 21, ticket 2374827, date 2010-11-20 (split across two seats) ---
$ws.Range("D20").Value = " '2010-11-20'"
$ws.Range("H20").Formula = "=CONCATENATE(`$M`$1,A20,`$O`$1,B20,`$O`$1,C20,`$O`$1,D20,`$N`$1)"

$ws.Range("A21").Value = 2374827
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 471
$ws.Range("D21").Value = " '2010-11-20'"
$ws.Range("H21").Formula = "=CONCATENATE(`$M`$1,A21,`$O`$1,B21,`$O`$1,C21,`$O`$1,D21,`$N`$1)"

# --- Group 7: row 23, ticket 2321423, date 2007-08-07 ---
$ws.Range("D23").Value = " '2007-08-07'"
$ws.Range("H23").Formula = "=CONCATENATE(`$M`$1,A23,`$O`$1,B23,`$O`$1,C23,`$O`$1,D23,`$N`$1)"

# 5. Column widths: col H keeps the old "bestFit" width; col M needs re-measuring since its
#    text grew longer.
$ws.Columns("M").AutoFit()

# 6. Selection moved to E11 in the saved file
$ws.Range("E11").Select()

Write-Host "done"
